$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Florida -- Miami-Dade County
$ws.Range("B14").Value = 44039
$ws.Range("C14").Value = 107315
$ws.Range("D14").Value = 1404
$ws.Range("E14").Value = 9879
$ws.Range("G14").Value = 18.27
$ws.Range("K14").Value = 54066

# Row 15: Florida -- Orange County
$ws.Range("B15").Value = 44039
$ws.Range("C15").Value = 27393
$ws.Range("D15").Value = 174
$ws.Range("E15").Value = 3513
$ws.Range("G15").Value = 25.12
$ws.Range("K15").Value = 13983

# Row 18: Maryland - was all empty/error, now filled in with results
$ws.Range("B18").Value = 44039
$ws.Range("B18").NumberFormat = "YYYY-MM-DD"
$ws.Range("C18").Value = 84876
$ws.Range("D18").Value = 3315
$ws.Range("E18").Value = 25440
$ws.Range("F18").Value = 1357
$ws.Range("G18").Value = 36.16
$ws.Range("H18").Value = 41.08
$ws.Range("K18").Value = 70362
$ws.Range("L18").Value = 3303
$ws.Range("O18").Value = "Success!"

# Row 24: Vermont
$ws.Range("B24").Value = 44039
$ws.Range("C24").Value = 1402
$ws.Range("G24").Value = 11.67
$ws.Range("K24").Value = 1362
